$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("baseline")

# New row of results for "deberta + cnn"
$ws.Range("B10").Value = "deberta + cnn"
$ws.Range("C10").Value = 1.7863
$ws.Range("D10").Value = 0.632
$ws.Range("E10").Value = 0.5866
$ws.Range("F10").Value = 0.632
$ws.Range("G10").Value = 0.5899
$ws.Range("H10").Value = 0.3124
$ws.Range("I10").Value = 0.3167
$ws.Range("J10").Value = 0.4711

# Match the number format used by the other metric columns
$ws.Range("C10:J10").NumberFormat = "0.0000"

# Update the active cell / selection on the sheet
$ws.Activate()
$ws.Range("G16").Select()
